$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that removes clientUpdateProhibited (originally row 22, "rem"/"clientUpdateProhibited"/
# "pass"/"EPP_UNEXPECTED_COMMAND_FAILURE") needs to move up to directly follow the row that adds
# clientUpdateProhibited (row 10), so it becomes row 11. This also means it loses its
# passOrFail/errorCode values (those cells become empty, inheriting from context as the sheet's
# notes describe), while all the "add" rows that used to be 11-21 shift down to 12-22.

# 1. Insert a new blank row at row 11, pushing the existing rows 11-34 down to 12-35.
$ws.Rows.Item(11).Insert()

# 2. Populate the new row 11 with the moved data (only action + status; passOrFail/errorCode left
#    completely blank/empty, matching the other rows in the sheet that don't set every column).
$ws.Range("B11").Value = "rem"
$ws.Range("C11").Value = "clientUpdateProhibited"
$ws.Range("D11:E11").Clear()

# 3. The original "rem clientUpdateProhibited" row is now duplicated at row 23 (old row 22 shifted
#    down by the insert). Remove that duplicate row, shifting rows 24-35 back up to 23-34.
$ws.Rows.Item(23).Delete()

# 4. Deleting a row out of the table's body can shrink the table definition by one row; put the
#    table back to its original B7:E34 extent now that the net row count is unchanged.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B7:E34"))

# 5. Update the selection to match the author's final cursor position.
$ws.Range("E17").Select()
